# "Generate Report for handoff"
#
# The handoff transform for the 1d7d6837-...-...md source file failed, so
# the localization-status report needs to reflect that:
#   - Status goes from "Ready for handoff" to "Handoff transform failed"
#     (shown both on the Overview sheet and on each per-language sheet)
#   - the per-language sheet's "Latest Handoff File" hyperlink/cell is
#     cleared out (no handoff file was produced)
#   - "Latest Handoff Datetime" resets to the zero-date sentinel
#   - "Handoff Reason" flips from "Include" to "Ignored"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove just the hyperlink anchored at C2 (Hyperlinks.Delete() on a
    # range deletes every hyperlink on the sheet, so find the matching one
    # explicitly instead).
    $linksToRemove = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $linksToRemove += $hl
        }
    }
    foreach ($hl in $linksToRemove) {
        $hl.Delete()
    }

    $ws.Range("B2").Value = "Handoff transform failed"
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
